$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "U" value in column D for the data rows (2-6), stored as a
# shared string (matches the new shared string entry "U" added to the
# workbook), and remove the now-unused sire/dam values in columns F and G.
$ws.Range("D2").Value = "U"
$ws.Range("D3").Value = "U"
$ws.Range("D4").Value = "U"
$ws.Range("D5").Value = "U"
$ws.Range("D6").Value = "U"

$ws.Range("F2:G6").ClearContents()

# Update the saved selection/active cell to D7.
$ws.Range("D7").Select()
